# Add a new "APR-2022" worksheet after "MAR-2022", populate it with the
# first three weeks of April-2022 tracker rows, and update the
# previously-last sheet's view state (it is no longer the active tab).

$wb = $excel.ActiveWorkbook

# --- 1) Fix up the view state of the (current) last sheet, MAR-2022,
#        BEFORE we touch anything else: selecting a cell does not disturb
#        row heights, so do it first while the sheet is still pristine. ---
$wsMar = $wb.Worksheets.Item("MAR-2022")
$wsMar.Activate()
$wsMar.Range("D37").Select()

# --- 2) Duplicate MAR-2022 to get an exact copy of all formatting
#        (column widths, cell styles, shared styles) placed right after it. ---
$wsMar.Copy([System.Reflection.Missing]::Value, $wsMar)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "APR-2022"

# --- 3) Overwrite header styles are already correct (copied). Update the
#        data rows with the April values. Columns: A=No, B=Date,
#        C=Application, D=Task, E=% completion, F=Status, G=Comments. ---

# Row 2 - Apr 1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 44652
$ws.Range("D2").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. "
$ws.Range("C2").Value = "Sanity Testing and Starshield app"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "Completed"
$ws.Range("G2").Value = $null

# Row 3 - Apr 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 44653
$ws.Range("D3").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. "
$ws.Range("C3").Value = "Sanity Testing and Starshield app"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Completed"
$ws.Range("G3").Value = $null

# Row 4 - Apr 3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 44654
$ws.Range("D4").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. "
$ws.Range("C4").Value = "Sanity Testing and Starshield app"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Completed"
$ws.Range("G4").Value = $null

# Row 7 (MJA texts) is written before row 5/6 (Contract Management texts) so
# that new shared strings get appended to xl/sharedStrings.xml in the same
# order as the target workbook: MJA-long, MJA-short, Contract-short,
# Contract-long.
$ws.Range("D7").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. `nRegression and retesting on MJA application."
$ws.Range("C7").Value = "Sanity Testing, Starshield app, and MJA"
$ws.Range("C5").Value = "Sanity Testing,  Starshield app and Contract Management"
$ws.Range("D5").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. `nRegression and retesting on Contract Management application."

# Row 5 - Apr 4 (Contract Management)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44655
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "Completed"
$ws.Range("G5").Value = $null

# Row 6 - Apr 5 (Contract Management)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44656
$ws.Range("C6").Value = "Sanity Testing,  Starshield app and Contract Management"
$ws.Range("D6").Value = "Sanity testing on B2C/B2B app, QMVAR site, GSS site and Hayaai site. `nRegression and retesting on Starshield app. `nRegression and retesting on Contract Management application."
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "Completed"
$ws.Range("G6").Value = $null

# Row 7 - Apr 6 (MJA)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44657
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "Completed"
$ws.Range("G7").Value = $null

# --- 4) Rows 8-21: only the sequence number / date are known; clear the
#        remaining (copied) text and reset their style to the plain
#        bordered style ("1"), matching unused tracker rows elsewhere in
#        the workbook. Borrow that style from G2 (already style "1"). ---
$ws.Range("G2").Copy()
$ws.Range("C8:F21").PasteSpecial(-4122)
$ws.Range("C8:F21").ClearContents()
$excel.CutCopyMode = 0

$dates = 44658,44659,44660,44661,44662,44663,44664,44665,44666,44667,44668,44669,44670,44671
$r = 8
$n = 7
foreach ($d in $dates) {
    $ws.Range("A" + $r).Value = $n
    $ws.Range("B" + $r).Value = $d
    $r = $r + 1
    $n = $n + 1
}

# --- 5) Drop the trailing rows copied from MAR-2022 (23-32) so the sheet
#        ends at row 22, then blank out row 22 entirely (same bordered
#        style, no values at all). ---
$ws.Range("A23:G32").EntireRow.Delete()

$ws.Range("G2").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$ws.Range("A22:G22").ClearContents()
$excel.CutCopyMode = 0

# --- 6) Column widths: match the narrower auto-fit widths the real
#        workbook ended up with for the Application/Task columns now
#        that their text is shorter. ---
$ws.Columns.Item(3).ColumnWidth = 17.307291666666668
$ws.Columns.Item(4).ColumnWidth = 51.166666666666664

# --- 7) Final view state for the new sheet: it becomes the active tab,
#        scrolled to the top, with G5 selected. ---
$ws.Activate()
$ws.Range("G5").Select()
